$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data as scraped on Fri Nov 24 09:31:25 UTC 2023

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '37.582.73'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.23%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.087.58'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.51%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '234.37'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.25%  '
$ws.Range("E6").Value = '  +0.77%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '58.20'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.11%  '
$ws.Range("E8").Value = '  +0.04%  '
$ws.Range("E9").Value = '  +1.12%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0781'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.37%  '
$ws.Range("E11").Value = '  +1.28%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '2.396.23'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.48%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '14.50'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.60%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '21.24'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.39%  '
$ws.Range("E15").Value = '  +0.62%  '
$ws.Range("E16").Value = '  +0.71%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.091.19'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.67%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '37.530.44'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.23%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.19'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.31%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '69.76'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.16%  '
$ws.Range("E21").Value = '  +0.62%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '226.80'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.14%  '
$ws.Range("E23").Value = '  +0.20%  '
$ws.Range("E24").Value = '  +2.94%  '
$ws.Range("E25").Value = '  -2.84%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '169.14'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.11%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.92'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.13%  '
$ws.Range("B28").Value = 'Kaspa'
$ws.Range("C28").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.133'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +5.01%  '
$ws.Range("B29").Value = 'ImmutableX'
$ws.Range("C29").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.44'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -4.55%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '19.29'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.24%  '
$ws.Range("E31").Value = '  +0.24%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.65'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.27%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0620'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.80%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.58'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.16%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.53'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.11%  '
$ws.Range("E36").Value = '  +5.19%  '
$ws.Range("E37").Value = '  +1.07%  '
$ws.Range("E38").Value = '  +0.09%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.59'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -5.00%  '
$ws.Range("E40").Value = '  -0.71%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0958'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.03%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.484.63'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.38%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '97.27'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.41%  '
$ws.Range("E44").Value = '  +0.17%  '
$ws.Range("E45").Value = '  -1.37%  '
$ws.Range("E46").Value = '  -10.07%  '
$ws.Range("E47").Value = '  +0.77%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '15.52'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.22%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.29'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.09%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.02'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.64%  '
$ws.Range("B51").Value = 'RocketPoolETH'
$ws.Range("C51").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.282.81'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.51%  '
